# Update localization status report:
#  - change status text from "Ready for handoff" to "In Translation"
#    on the Overview sheet (per-language status columns) and on each
#    per-language detail sheet (zh-cn, de-de).
#  - shrink the now-narrower "Status" columns to fit the new text
#    (mirrors Excel re-computing the column width for the edited cells).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newWidth = 13.4101845877511

# --- Overview sheet: columns E (zh-cn status) and F (de-de status) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet: column C (Status) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet: column C (Status) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
